# Instructions.xlsx: add a "Sheet2" with Data Type + control-bit columns
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# New worksheet, placed right after Sheet1
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Header row (A:Name, C:Opcode, D:Operation, E:Description, F-J: control bits)
$ws2.Cells.Item(1,1).Value = 'Name'
$ws2.Cells.Item(1,3).Value = 'Opcode'
$ws2.Cells.Item(1,4).Value = 'Operation'
$ws2.Cells.Item(1,5).Value = 'Description'
$ws2.Cells.Item(1,6).Value = 'RegWrite'
$ws2.Cells.Item(1,7).Value = 'MemtoReg'
$ws2.Cells.Item(1,8).Value = 'MemWrite'
$ws2.Cells.Item(1,9).Value = 'ALUSrc'
$ws2.Cells.Item(1,10).Value = 'SignExtend'
$ws2.Cells.Item(1,2).Value = 'Data Type'

# New "Data Type" column (B): R = register, D = displacement/memory, B = branch
$ws2.Cells.Item(3,2).Value = 'D'
$ws2.Cells.Item(2,2).Value = 'R'
$ws2.Cells.Item(4,2).Value = 'R'
$ws2.Cells.Item(5,2).Value = 'D'
$ws2.Cells.Item(6,2).Value = 'R'
$ws2.Cells.Item(7,2).Value = 'R'
$ws2.Cells.Item(8,2).Value = 'R'
$ws2.Cells.Item(9,2).Value = 'R'
$ws2.Cells.Item(10,2).Value = 'D'
$ws2.Cells.Item(11,2).Value = 'D'
$ws2.Cells.Item(12,2).Value = 'D'
$ws2.Cells.Item(13,2).Value = 'B'

# Name / Opcode / Operation / Description (copied from Sheet1) + new control-bit columns
$ws2.Cells.Item(2,1).Value = 'ADD'
$ws2.Cells.Item(2,3).Value = 0
$ws2.Cells.Item(2,4).Value = 'GPR[Rd] = GPR[Rs1] + GPR[Rs2]'
$ws2.Cells.Item(2,5).Value = 'Add'
$ws2.Cells.Item(2,6).Value = 1
$ws2.Cells.Item(2,7).Value = 0
$ws2.Cells.Item(2,8).Value = 0
$ws2.Cells.Item(2,9).Value = 0
$ws2.Cells.Item(2,10).Value = 0
$ws2.Cells.Item(3,1).Value = 'ADDM'
$ws2.Cells.Item(3,3).Value = 1
$ws2.Cells.Item(3,4).Value = 'GPR[Rd] = GPR[Rd] + MM[PC + Short_Offset]'
$ws2.Cells.Item(3,5).Value = 'Add from memory'
$ws2.Cells.Item(4,1).Value = 'SUB'
$ws2.Cells.Item(4,3).Value = 2
$ws2.Cells.Item(4,4).Value = 'GPR[Rd] = GPR[Rs1] - GPR[Rs2]'
$ws2.Cells.Item(4,5).Value = 'Subtract'
$ws2.Cells.Item(5,1).Value = 'SUBM'
$ws2.Cells.Item(5,3).Value = 3
$ws2.Cells.Item(5,4).Value = 'GPR[Rd] = GPR[Rd] - MM[PC + Short_Offset]'
$ws2.Cells.Item(5,5).Value = 'Subtract from memory'
$ws2.Cells.Item(6,1).Value = 'AND'
$ws2.Cells.Item(6,3).Value = 4
$ws2.Cells.Item(6,4).Value = 'GPR[Rd] = GPR[Rs1] and GPR[Rs2]'
$ws2.Cells.Item(6,5).Value = 'And'
$ws2.Cells.Item(7,1).Value = 'SHL'
$ws2.Cells.Item(7,3).Value = 5
$ws2.Cells.Item(7,4).Value = 'GPR[Rd] = shift_left(GPR[Rs1]) by GPR[Rs2]_3-0'
$ws2.Cells.Item(7,5).Value = 'Shift Left'
$ws2.Cells.Item(8,1).Value = 'SHRA'
$ws2.Cells.Item(8,3).Value = 6
$ws2.Cells.Item(8,4).Value = 'GPR[Rd] = shift_right(GPR[Rs1]) by GPR[Rs2]_3-0'
$ws2.Cells.Item(8,5).Value = 'Shift Right'
$ws2.Cells.Item(9,1).Value = 'OR'
$ws2.Cells.Item(9,3).Value = 7
$ws2.Cells.Item(9,4).Value = 'GPR[Rd] = GPR[Rs1] or GPR[Rs2]'
$ws2.Cells.Item(9,5).Value = 'Or'
$ws2.Cells.Item(10,1).Value = 'NOT'
$ws2.Cells.Item(10,3).Value = 8
$ws2.Cells.Item(10,4).Value = 'GPR[Rd] = not MM[PC + Short_Offset]'
$ws2.Cells.Item(10,5).Value = 'Not'
$ws2.Cells.Item(11,1).Value = 'LD'
$ws2.Cells.Item(11,3).Value = 9
$ws2.Cells.Item(11,4).Value = 'GPR[Rd] = MM[PC + Short_Offset]'
$ws2.Cells.Item(11,5).Value = 'Load from memory'
$ws2.Cells.Item(11,6).Value = 1
$ws2.Cells.Item(11,7).Value = 1
$ws2.Cells.Item(11,8).Value = 0
$ws2.Cells.Item(12,1).Value = 'ST'
$ws2.Cells.Item(12,3).Value = 10
$ws2.Cells.Item(12,4).Value = 'MM[PC + Short_Offset] = GPR[Rd]'
$ws2.Cells.Item(12,5).Value = 'Store to memory'
$ws2.Cells.Item(12,6).Value = 0
$ws2.Cells.Item(12,7).Value = 0
$ws2.Cells.Item(12,8).Value = 1
$ws2.Cells.Item(13,1).Value = 'BR'
$ws2.Cells.Item(13,3).Value = 11
$ws2.Cells.Item(13,4).Value = 'PC = PC + Long_Offset'
$ws2.Cells.Item(13,5).Value = 'Branch'
$ws2.Cells.Item(14,1).Value = 'JSR'
$ws2.Cells.Item(14,3).Value = 12
$ws2.Cells.Item(14,4).Value = 'GPR[Rd] = PC; PC = PC + Short_Offset'
$ws2.Cells.Item(15,1).Value = 'RTS'
$ws2.Cells.Item(15,3).Value = 13
$ws2.Cells.Item(15,4).Value = 'PC = GPR[Rd] + Short_Offset'
$ws2.Cells.Item(16,1).Value = 'CLK'
$ws2.Cells.Item(16,3).Value = 14
$ws2.Cells.Item(16,4).Value = 'Set timer to MM[PC + Long_Offset]'
$ws2.Cells.Item(17,1).Value = 'LPSW'
$ws2.Cells.Item(17,3).Value = 15
$ws2.Cells.Item(17,4).Value = 'PSW = MM[PC + Long_Offset]'

# Turn the range into a table, same way Table1 backs Sheet1
$tbl2 = $ws2.ListObjects.Add(1, $ws2.Range("A1:J17"), $null, 1)
$tbl2.Name = "Table13"
$tbl2.TableStyle = "TableStyleLight1"

# Sheet1 no longer shows the stray C20 selection - select the whole table range instead
[void]$ws1.Range("A1:D17").Select()

# Sheet2 is the tab that ends up active/selected
$ws2.Activate()
[void]$ws2.Range("B3").Select()
